$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (column B); the columns to its
# right (Not living with children, Living with children, Not known /
# missing, Total) shift left to take its place.
$ws.Columns("B").Delete()
